$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Range("H32").Value = 346.1875
$ws.Range("I32").Value = 316.83334
$ws.Range("J32").Value = 363.8
$ws.Range("K32").Value = 316.83334
$ws.Range("L32").Value = 363.8
$ws.Range("M32").Value = 9.166659999999979
$ws.Range("N32").Value = -1015.8
# Row 98
$ws.Range("H98").Value = 857.7143
$ws.Range("I98").Value = 857.7143
$ws.Range("K98").Value = 857.7143
$ws.Range("M98").Value = 640.2857
# Row 107
$ws.Range("H107").Value = 366.42105
$ws.Range("I107").Value = 418.9375
$ws.Range("J107").Value = 86.333336
$ws.Range("K107").Value = 418.9375
$ws.Range("L107").Value = 86.333336
$ws.Range("M107").Value = 1501.0625
$ws.Range("N107").Value = -3926.333336
# Row 122
$ws.Range("H122").Value = 857.7143
$ws.Range("I122").Value = 857.7143
$ws.Range("K122").Value = 2573.1429
$ws.Range("M122").Value = -123.1428999999998
# Row 129
$ws.Range("H129").Value = 2371.8333
$ws.Range("I129").Value = 8270.77
$ws.Range("J129").Value = 924.9245
$ws.Range("K129").Value = 24812.31
$ws.Range("L129").Value = 2774.7735
$ws.Range("M129").Value = -19812.31
$ws.Range("N129").Value = -12774.7735
# Row 137
$ws.Range("H137").Value = 1792.9231
$ws.Range("I137").Value = 1889.7142
$ws.Range("J137").Value = 1680
$ws.Range("K137").Value = 5669.142599999999
$ws.Range("L137").Value = 5040
$ws.Range("M137").Value = -3119.142599999999
$ws.Range("N137").Value = -10140
# Row 138
$ws.Range("H138").Value = 3609.8276
$ws.Range("I138").Value = 1996.6875
$ws.Range("J138").Value = 5595.231
$ws.Range("K138").Value = 5990.0625
$ws.Range("L138").Value = 16785.693
$ws.Range("M138").Value = -850.0625
$ws.Range("N138").Value = -27065.693

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 102
$ws.Range("H102").Value = 113995.11
$ws.Range("I102").Value = 252244.75
$ws.Range("J102").Value = 3395.4
$ws.Range("K102").Value = 252244.75
$ws.Range("L102").Value = 3395.4
$ws.Range("M102").Value = -250622.75
$ws.Range("N102").Value = -6639.4
# Row 122
$ws.Range("H122").Value = 1625.0476
$ws.Range("I122").Value = 1696.1428
$ws.Range("J122").Value = 1482.8572
$ws.Range("K122").Value = 5088.428400000001
$ws.Range("L122").Value = 4448.571599999999
$ws.Range("M122").Value = -2638.428400000001
$ws.Range("N122").Value = -9348.571599999999
# Row 132
$ws.Range("H132").Value = 4247.7144
$ws.Range("I132").Value = 4473.579
$ws.Range("K132").Value = 13420.737
$ws.Range("M132").Value = -10890.737
# Row 134
$ws.Range("H134").Value = 37211.25
$ws.Range("J134").Value = 37211.25
$ws.Range("L134").Value = 37211.25
$ws.Range("N134").Value = -47351.25
# Row 135
$ws.Range("H135").Value = 42049.8
$ws.Range("J135").Value = 42049.8
$ws.Range("L135").Value = 42049.8
$ws.Range("N135").Value = -52189.8
# Row 137
$ws.Range("H137").Value = 47000
$ws.Range("J137").Value = 47000
$ws.Range("L137").Value = 47000
$ws.Range("N137").Value = -57200

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 71969.81
$ws.Range("I86").Value = 126390.22
$ws.Range("J86").Value = 2000.7142
$ws.Range("K86").Value = 126390.22
$ws.Range("L86").Value = 2000.7142
$ws.Range("M86").Value = -125267.22
$ws.Range("N86").Value = -4246.7142
# Row 89
$ws.Range("H89").Value = 71969.81
$ws.Range("I89").Value = 126390.22
$ws.Range("J89").Value = 2000.7142
$ws.Range("K89").Value = 631951.1
$ws.Range("L89").Value = 10003.571
$ws.Range("M89").Value = -626335.1
$ws.Range("N89").Value = -21235.571
# Row 107
$ws.Range("H107").Value = 41667304
$ws.Range("I107").Value = 76923304
$ws.Range("J107").Value = 1127.2727
$ws.Range("K107").Value = 76923304
$ws.Range("L107").Value = 1127.2727
$ws.Range("M107").Value = -76921384
$ws.Range("N107").Value = -4967.2727
# Row 133
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 105
$ws.Range("H105").Value = 1705.3334
$ws.Range("I105").Value = 1775.8
$ws.Range("K105").Value = 1775.8
$ws.Range("M105").Value = -28.79999999999995
# Row 138
$ws.Range("H138").Value = 61482.418
$ws.Range("J138").Value = 61482.418
$ws.Range("L138").Value = 61482.418
$ws.Range("N138").Value = -71762.41800000001
# Row 140
$ws.Range("H140").Value = 59500
$ws.Range("J140").Value = 59500
$ws.Range("L140").Value = 59500
$ws.Range("N140").Value = -69860

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 37
$ws.Range("H37").Value = 610740.5600000001
$ws.Range("J37").Value = 610740.5600000001
$ws.Range("L37").Value = 1832221.68
$ws.Range("N37").Value = -1832445.68
# Row 38
$ws.Range("H38").Value = 34.3125
$ws.Range("J38").Value = 50.57143
$ws.Range("L38").Value = 151.71429
$ws.Range("N38").Value = -845.71429
# Row 58
$ws.Range("H58").Value = 2266.6667
$ws.Range("J58").Value = 2950
$ws.Range("L58").Value = 8850
$ws.Range("N58").Value = -9106
# Row 131
$ws.Range("H131").Value = 6356.18
$ws.Range("I131").Value = 1211.6666
$ws.Range("J131").Value = 6684.553
$ws.Range("K131").Value = 3634.9998
$ws.Range("L131").Value = 20053.659
$ws.Range("M131").Value = 1405.0002
$ws.Range("N131").Value = -30133.659
# Row 132
$ws.Range("H132").Value = 2039.9231
$ws.Range("I132").Value = 682.63635
$ws.Range("K132").Value = 6143.72715
$ws.Range("M132").Value = -3613.72715
# Row 138
$ws.Range("H138").Value = 11235.363
$ws.Range("I138").Value = 16012.857
$ws.Range("J138").Value = 2874.75
$ws.Range("K138").Value = 48038.571
$ws.Range("L138").Value = 8624.25
$ws.Range("M138").Value = -42898.571
$ws.Range("N138").Value = -18904.25
# Row 139
$ws.Range("H139").Value = 2373.7058
$ws.Range("I139").Value = 1070
$ws.Range("J139").Value = 3532.5557
$ws.Range("K139").Value = 3210
$ws.Range("L139").Value = 10597.6671
$ws.Range("M139").Value = 1930
$ws.Range("N139").Value = -20877.6671
# Row 140
$ws.Range("H140").Value = 4873.207
$ws.Range("I140").Value = 6510.5557
$ws.Range("J140").Value = 2193.9092
$ws.Range("K140").Value = 19531.6671
$ws.Range("L140").Value = 6581.7276
$ws.Range("M140").Value = -14351.6671
$ws.Range("N140").Value = -16941.7276
# Row 141
$ws.Range("H141").Value = 9145.643
$ws.Range("I141").Value = 10280.818
$ws.Range("J141").Value = 4983.3335
$ws.Range("K141").Value = 30842.454
$ws.Range("L141").Value = 14950.0005
$ws.Range("M141").Value = -25662.454
$ws.Range("N141").Value = -25310.0005

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 26
$ws.Range("H26").Value = 8570
$ws.Range("J26").Value = 8570
$ws.Range("L26").Value = 8570
$ws.Range("N26").Value = -9130
# Row 50
$ws.Range("H50").Value = 8570
$ws.Range("J50").Value = 8570
$ws.Range("L50").Value = 8570
$ws.Range("N50").Value = -9566
# Row 52
$ws.Range("H52").Value = 11622.857
$ws.Range("I52").Value = 4130
$ws.Range("J52").Value = 14620
$ws.Range("K52").Value = 4130
$ws.Range("L52").Value = 14620
$ws.Range("M52").Value = -3871
$ws.Range("N52").Value = -15138
# Row 126
$ws.Range("H126").Value = 3118.8
$ws.Range("I126").Value = 4990
$ws.Range("J126").Value = 1871.3334
$ws.Range("K126").Value = 14970
$ws.Range("L126").Value = 5614.0002
$ws.Range("M126").Value = -12500
$ws.Range("N126").Value = -10554.0002
# Row 133
$ws.Range("H133").Value = 49666.668
$ws.Range("J133").Value = 49666.668
$ws.Range("L133").Value = 49666.668
$ws.Range("N133").Value = -59786.668
# Row 135
$ws.Range("H135").Value = 33652.6
$ws.Range("J135").Value = 33652.6
$ws.Range("L135").Value = 33652.6
$ws.Range("N135").Value = -43792.6

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 43392.082
$ws.Range("I40").Value = 101148.7
$ws.Range("J40").Value = 2137.3572
$ws.Range("K40").Value = 101148.7
$ws.Range("L40").Value = 2137.3572
$ws.Range("M40").Value = -101012.7
$ws.Range("N40").Value = -2409.3572
# Row 61
$ws.Range("H61").Value = 1784.55
$ws.Range("I61").Value = 1659
$ws.Range("J61").Value = 1972.875
$ws.Range("K61").Value = 1659
$ws.Range("L61").Value = 1972.875
$ws.Range("M61").Value = -1457
$ws.Range("N61").Value = -2376.875
# Row 113
$ws.Range("H113").Value = 1784.55
$ws.Range("I113").Value = 1659
$ws.Range("J113").Value = 1972.875
$ws.Range("K113").Value = 1659
$ws.Range("L113").Value = 1972.875
$ws.Range("M113").Value = 511
$ws.Range("N113").Value = -6312.875
# Row 122
$ws.Range("H122").Value = 2570.2942
$ws.Range("I122").Value = 2490.8333
$ws.Range("J122").Value = 2761
$ws.Range("K122").Value = 7472.499899999999
$ws.Range("L122").Value = 8283
$ws.Range("M122").Value = -5022.499899999999
$ws.Range("N122").Value = -13183
# Row 132
$ws.Range("H132").Value = 4390
$ws.Range("I132").Value = 4340.6313
$ws.Range("K132").Value = 13021.8939
$ws.Range("M132").Value = -10491.8939
# Row 133
$ws.Range("H133").Value = 33339.3
$ws.Range("J133").Value = 33339.3
$ws.Range("L133").Value = 33339.3
$ws.Range("N133").Value = -38399.3
# Row 134
$ws.Range("H134").Value = 62132.715
$ws.Range("J134").Value = 62132.715
$ws.Range("L134").Value = 62132.715
$ws.Range("N134").Value = -72272.715
# Row 136
$ws.Range("H136").Value = 1907
$ws.Range("I136").Value = 1842.0769
$ws.Range("K136").Value = 5526.2307
$ws.Range("M136").Value = -2976.2307
# Row 137
$ws.Range("H137").Value = 28800
$ws.Range("J137").Value = 35750
$ws.Range("L137").Value = 35750
$ws.Range("N137").Value = -45950
# Row 139
$ws.Range("H139").Value = 65396.668
$ws.Range("J139").Value = 65396.668
$ws.Range("L139").Value = 65396.668
$ws.Range("N139").Value = -75676.66800000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 472
$ws.Range("I107").Value = 301.5
$ws.Range("J107").Value = 715.5714
$ws.Range("K107").Value = 904.5
$ws.Range("L107").Value = 2146.7142
$ws.Range("M107").Value = 1015.5
$ws.Range("N107").Value = -5986.7142
# Row 113
$ws.Range("H113").Value = 754.9091
$ws.Range("I113").Value = 568.1818
$ws.Range("J113").Value = 941.63635
$ws.Range("K113").Value = 1704.5454
$ws.Range("L113").Value = 2824.90905
$ws.Range("M113").Value = 465.4546
$ws.Range("N113").Value = -7164.90905
# Row 140
$ws.Range("H140").Value = 61532.855
$ws.Range("J140").Value = 61532.855
$ws.Range("L140").Value = 61532.855
$ws.Range("N140").Value = -71892.85500000001
